$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue $ws.Cells.Item(2, 4) '37.757.71'
Set-TextValue $ws.Cells.Item(2, 5) '  +1.38%  '
Set-TextValue $ws.Cells.Item(3, 4) '2.086.61'
Set-TextValue $ws.Cells.Item(3, 5) '  +1.20%  '
Set-TextValue $ws.Cells.Item(5, 4) '232.53'
Set-TextValue $ws.Cells.Item(5, 5) '  -0.66%  '
Set-TextValue $ws.Cells.Item(6, 4) '0.623'
Set-TextValue $ws.Cells.Item(6, 5) '  -0.38%  '
Set-TextValue $ws.Cells.Item(7, 5) '  -0.03%  '
Set-TextValue $ws.Cells.Item(8, 4) '57.26'
Set-TextValue $ws.Cells.Item(8, 5) '  +1.22%  '
Set-TextValue $ws.Cells.Item(9, 5) '  +1.61%  '
Set-TextValue $ws.Cells.Item(10, 4) '0.0777'
Set-TextValue $ws.Cells.Item(10, 5) '  +2.23%  '
Set-TextValue $ws.Cells.Item(11, 5) '  +2.97%  '
Set-TextValue $ws.Cells.Item(12, 4) '2.384.17'
Set-TextValue $ws.Cells.Item(12, 5) '  +0.78%  '
Set-TextValue $ws.Cells.Item(13, 5) '  -1.32%  '
Set-TextValue $ws.Cells.Item(14, 4) '21.06'
Set-TextValue $ws.Cells.Item(14, 5) '  +2.48%  '
Set-TextValue $ws.Cells.Item(15, 4) '0.765'
Set-TextValue $ws.Cells.Item(15, 5) '  -1.21%  '
Set-TextValue $ws.Cells.Item(16, 4) '5.22'
Set-TextValue $ws.Cells.Item(16, 5) '  +2.24%  '
Set-TextValue $ws.Cells.Item(17, 4) '2.076.74'
Set-TextValue $ws.Cells.Item(17, 5) '  +0.71%  '
Set-TextValue $ws.Cells.Item(18, 4) '37.669.84'
Set-TextValue $ws.Cells.Item(18, 5) '  +1.15%  '
Set-TextValue $ws.Cells.Item(19, 4) '6.12'
Set-TextValue $ws.Cells.Item(19, 5) '  -3.34%  '
Set-TextValue $ws.Cells.Item(20, 4) '70.70'
Set-TextValue $ws.Cells.Item(20, 5) '  +1.88%  '
Set-TextValue $ws.Cells.Item(21, 4) '0.0₃0820'
Set-TextValue $ws.Cells.Item(21, 5) '  +1.46%  '
Set-TextValue $ws.Cells.Item(22, 5) '  +0.81%  '
Set-TextValue $ws.Cells.Item(23, 5) '  +0.00%  '
Set-TextValue $ws.Cells.Item(24, 5) '  -1.59%  '
Set-TextValue $ws.Cells.Item(25, 4) '2.38'
Set-TextValue $ws.Cells.Item(25, 5) '  -0.31%  '
Set-TextValue $ws.Cells.Item(26, 4) '168.16'
Set-TextValue $ws.Cells.Item(26, 5) '  +1.29%  '
Set-TextValue $ws.Cells.Item(27, 4) '0.141'
Set-TextValue $ws.Cells.Item(27, 5) '  +11.98%  '
Set-TextValue $ws.Cells.Item(28, 4) '8.91'
Set-TextValue $ws.Cells.Item(28, 5) '  +1.90%  '
Set-TextValue $ws.Cells.Item(29, 5) '  -0.97%  '
Set-TextValue $ws.Cells.Item(30, 4) '19.44'
Set-TextValue $ws.Cells.Item(30, 5) '  +2.40%  '
Set-TextValue $ws.Cells.Item(31, 4) '0.118'
Set-TextValue $ws.Cells.Item(31, 5) '  +1.12%  '
Set-TextValue $ws.Cells.Item(32, 5) '  +3.76%  '
Set-TextValue $ws.Cells.Item(33, 5) '  +1.43%  '
Set-TextValue $ws.Cells.Item(34, 4) '4.56'
Set-TextValue $ws.Cells.Item(34, 5) '  +0.77%  '
Set-TextValue $ws.Cells.Item(35, 4) '2.49'
Set-TextValue $ws.Cells.Item(35, 5) '  -0.02%  '
Set-TextValue $ws.Cells.Item(36, 5) '  +3.75%  '
Set-TextValue $ws.Cells.Item(37, 5) '  +4.66%  '
Set-TextValue $ws.Cells.Item(38, 5) '  +0.00%  '
Set-TextValue $ws.Cells.Item(39, 4) '5.40'
Set-TextValue $ws.Cells.Item(39, 5) '  -4.87%  '
Set-TextValue $ws.Cells.Item(40, 5) '  +6.32%  '
Set-TextValue $ws.Cells.Item(41, 4) '2.93'
Set-TextValue $ws.Cells.Item(41, 5) '  -0.50%  '
Set-TextValue $ws.Cells.Item(42, 4) '97.52'
Set-TextValue $ws.Cells.Item(42, 5) '  +1.68%  '
Set-TextValue $ws.Cells.Item(43, 5) '  +0.72%  '
Set-TextValue $ws.Cells.Item(44, 4) '1.451.49'
Set-TextValue $ws.Cells.Item(44, 5) '  -0.63%  '
Set-TextValue $ws.Cells.Item(45, 5) '  -0.19%  '
Set-TextValue $ws.Cells.Item(46, 5) '  +3.74%  '
Set-TextValue $ws.Cells.Item(47, 4) '4.06'
Set-TextValue $ws.Cells.Item(47, 5) '  -3.44%  '
Set-TextValue $ws.Cells.Item(48, 4) '15.58'
Set-TextValue $ws.Cells.Item(48, 5) '  +3.89%  '
Set-TextValue $ws.Cells.Item(49, 5) '  +2.61%  '
Set-TextValue $ws.Cells.Item(50, 4) '3.01'
Set-TextValue $ws.Cells.Item(50, 5) '  +1.93%  '
Set-TextValue $ws.Cells.Item(51, 4) '2.279.60'
Set-TextValue $ws.Cells.Item(51, 5) '  +1.20%  '
